$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# values (e.g. "551.75") are not auto-converted to numbers by Excel,
# matching the source data which stores these as literal strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '64.151.22'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '3.319.52'
$ws.Range('E3').Value = '  -0.31%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '551.75'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').Value = '172.60'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').Value = '0.616'
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').Value = '3.309.82'
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('E10').Value = '  +6.01%  '
$ws.Range('D11').Value = '0.628'
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('D12').Value = '53.31'
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').Value = '0.0000280'
$ws.Range('E13').Value = '  +3.36%  '
$ws.Range('D14').Value = '9.01'
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('D15').Value = '3.841.09'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('E16').Value = '  +2.70%  '
$ws.Range('D17').Value = '18.05'
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').Value = '3.331.12'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').Value = '64.027.66'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').Value = '11.63'
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('D21').Value = '0.980'
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('D22').Value = '455.24'
$ws.Range('E22').Value = '  +6.65%  '
$ws.Range('D23').Value = '5.00'
$ws.Range('E23').Value = '  +7.72%  '
$ws.Range('D24').Value = '4.05'
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('D25').Value = '86.44'
$ws.Range('E25').Value = '  +2.81%  '
$ws.Range('D26').Value = '13.72'
$ws.Range('E26').Value = '  +1.68%  '
$ws.Range('E27').Value = '  +1.14%  '
$ws.Range('D28').Value = '10.61'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('D29').Value = '8.56'
$ws.Range('E29').Value = '  +0.19%  '
$ws.Range('D30').Value = '30.76'
$ws.Range('E30').Value = '  +3.90%  '
$ws.Range('D31').Value = '6.51'
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('D32').Value = '11.39'
$ws.Range('E32').Value = '  +0.09%  '
$ws.Range('D33').Value = '61.70'
$ws.Range('E33').Value = '  +6.25%  '
$ws.Range('D34').Value = '563.92'
$ws.Range('E34').Value = '  -5.00%  '
$ws.Range('E35').Value = '  -0.26%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = '0.141'
$ws.Range('E37').Value = '  -1.71%  '
$ws.Range('D38').Value = '3.50'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').Value = '35.09'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('D40').Value = '0.364'
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('D41').Value = '0.0₃0725'
$ws.Range('E41').Value = '  -2.73%  '
$ws.Range('D42').Value = '3.037.71'
$ws.Range('E42').Value = '  -1.92%  '
$ws.Range('D43').Value = '0.0412'
$ws.Range('E43').Value = '  +1.75%  '
$ws.Range('D44').Value = '2.73'
$ws.Range('E44').Value = '  -1.28%  '
$ws.Range('D45').Value = '3.20'
$ws.Range('E45').Value = '  +0.96%  '
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('D48').Value = '0.999'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').Value = '140.99'
$ws.Range('E49').Value = '  +6.42%  '
$ws.Range('E50').Value = '  -3.19%  '
$ws.Range('D51').Value = '8.09'
